# Automatische test-sync: 2025-07-22 18:14:50
# Adds a new log entry (row 23) to the "Logs" sheet and refreshes the
# "Dashboard" category-count summary accordingly.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append new row to the Logs sheet -------------------------------------
$newRow = 23

$logs.Cells.Item($newRow, 1).Value = "Wat zijn jullie openingstijden?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Cells.Item($newRow, 4).Value = "Openingstijden / Locatie"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value = "2025-07-22 18:14:29"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Ja"

# Writing a multi-line value auto-sets a custom row height in this engine;
# AutoFit() clears that override so the row keeps the sheet's default
# height (matching the other rows, none of which carry an explicit "ht").
$logs.Rows.Item($newRow).AutoFit()

# --- Extend conditional formatting ranges to include the new row ----------
# Re-point every existing rule's AppliesTo range by one row (keeps the
# original rule order / priority / dxfId / formula intact).
function Extend-FormatConditionRange([string]$oldRange, [string]$newRange) {
    $fcs = $logs.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

Extend-FormatConditionRange "D2:D22" "D2:D23"
Extend-FormatConditionRange "G2:G22" "G2:G23"
Extend-FormatConditionRange "H2:H22" "H2:H23"
Extend-FormatConditionRange "I2:I22" "I2:I23"
Extend-FormatConditionRange "J2:J22" "J2:J23"

# --- Refresh the Dashboard category counts ---------------------------------
# The new mail is categorised as "Openingstijden / Locatie", so that
# category's count moves from 4 to 5 (tying it with "Productinformatie"),
# and the two rows swap places in the (count-sorted) summary table.
$dash.Cells.Item(2, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(2, 2).Value = 5

$dash.Cells.Item(4, 1).Value = "Productinformatie"
$dash.Cells.Item(4, 2).Value = 5
